$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102
$ws.Range("A102").Value = 130964526
$ws.Range("B102").Value = 79243
$ws.Range("D102").Value = 'NT'
$ws.Range("E102").Value = 6425
$ws.Range("F102").Value = 'Garnlav'
$ws.Range("G102").Value = 'Alectoria sarmentosa'
$ws.Range("H102").Value = '(Ach.) Ach.'
$ws.Range("Q102").Value = 509610
$ws.Range("R102").Value = 6719050
$ws.Range("AC102").Value = 'Enstaka . inventering åt vasa vind'

# Row 103
$ws.Range("A103").Value = 130964547
$ws.Range("B103").Value = 57881
$ws.Range("D103").Value = 'NT'
$ws.Range("E103").Value = 100049
$ws.Range("F103").Value = 'Spillkråka'
$ws.Range("G103").Value = 'Dryocopus martius'
$ws.Range("H103").Value = '(Linnaeus, 1758)'
$ws.Range("Q103").Value = 509495
$ws.Range("R103").Value = 6718877
$ws.Range("AC103").Value = 'Födosökspår . inventering åt vasa vind'

# Row 104
$ws.Range("A104").Value = 130964537
$ws.Range("B104").Value = 79243
$ws.Range("D104").Value = 'NT'
$ws.Range("E104").Value = 6425
$ws.Range("F104").Value = 'Garnlav'
$ws.Range("G104").Value = 'Alectoria sarmentosa'
$ws.Range("H104").Value = '(Ach.) Ach.'
$ws.Range("Q104").Value = 509822
$ws.Range("R104").Value = 6718960
$ws.Range("AC104").Value = 'Rikligt . inventering åt vasa vind'

# Row 106
$ws.Range("A106").Value = 130964541
$ws.Range("B106").Value = 91808
$ws.Range("D106").Value = 'NT'
$ws.Range("E106").Value = 1202
$ws.Range("F106").Value = 'Ullticka'
$ws.Range("G106").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H106").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q106").Value = 509703
$ws.Range("R106").Value = 6719018
$ws.Range("AC106").Value = 'Enstaka . inventering åt vasa vind'

# Row 105
$ws.Range("A105").Value = 130964642
$ws.Range("B105").Value = 99036
$ws.Range("D105").Value = 'LC'
$ws.Range("E105").Value = 221952
$ws.Range("F105").Value = 'Spindelblomster'
$ws.Range("G105").Value = 'Neottia cordata'
$ws.Range("H105").Value = '(L.) Rich.'
$ws.Range("Q105").Value = 509917
$ws.Range("R105").Value = 6719042
$ws.Range("AC105").Value = 'Måttliga förekomster . inventering åt vasa vind'

# Row 107
$ws.Range("A107").Value = 130964545
$ws.Range("B107").Value = 57073
$ws.Range("D107").Value = 'LC'
$ws.Range("E107").Value = 100138
$ws.Range("F107").Value = 'Tjäder'
$ws.Range("G107").Value = 'Tetrao urogallus'
$ws.Range("H107").Value = 'Linnaeus, 1758'
$ws.Range("Q107").Value = 509535
$ws.Range("R107").Value = 6718925
$ws.Range("AC107").Value = 'Spillning . inventering åt vasa vind'

# Row 108
$ws.Range("A108").Value = 130964647
$ws.Range("B108").Value = 92106
$ws.Range("D108").Value = 'NT'
$ws.Range("E108").Value = 658
$ws.Range("F108").Value = 'Rosenticka'
$ws.Range("G108").Value = 'Fomitopsis rosea'
$ws.Range("H108").Value = '(Alb. & Schwein.:Fr.) P.Karst.'
$ws.Range("Q108").Value = 509741
$ws.Range("R108").Value = 6718998
$ws.Range("AC108").Value = 'Måttliga förekomster . inventering åt vasa vind'

# Row 109
$ws.Range("A109").Value = 130964574
$ws.Range("B109").Value = 79243
$ws.Range("D109").Value = 'NT'
$ws.Range("E109").Value = 6425
$ws.Range("F109").Value = 'Garnlav'
$ws.Range("G109").Value = 'Alectoria sarmentosa'
$ws.Range("H109").Value = '(Ach.) Ach.'
$ws.Range("Q109").Value = 509667
$ws.Range("R109").Value = 6719184
$ws.Range("AC109").Value = 'Måttliga förekomster . inventering åt vasa vind'

# Row 110
$ws.Range("A110").Value = 130964544
$ws.Range("B110").Value = 57073
$ws.Range("D110").Value = 'LC'
$ws.Range("E110").Value = 100138
$ws.Range("F110").Value = 'Tjäder'
$ws.Range("G110").Value = 'Tetrao urogallus'
$ws.Range("H110").Value = 'Linnaeus, 1758'
$ws.Range("Q110").Value = 509543
$ws.Range("R110").Value = 6718926
$ws.Range("AC110").Value = 'Spillning . inventering åt vasa vind'

# Row 119
$ws.Range("A119").Value = 130964648
$ws.Range("B119").Value = 92267
$ws.Range("D119").Value = 'VU'
$ws.Range("E119").Value = 1209
$ws.Range("F119").Value = 'Rynkskinn'
$ws.Range("G119").Value = 'Hermanssonia centrifuga'
$ws.Range("H119").Value = '(P. Karst.) Zmitr.'
$ws.Range("Q119").Value = 509744
$ws.Range("R119").Value = 6718982
$ws.Range("AC119").Value = 'Måttliga förekomster . inventering åt vasa vind'

# Row 120
$ws.Range("A120").Value = 130964649
$ws.Range("B120").Value = 98930
$ws.Range("D120").Value = 'LC'
$ws.Range("E120").Value = 219790
$ws.Range("F120").Value = 'Fläcknycklar'
$ws.Range("G120").Value = 'Dactylorhiza maculata'
$ws.Range("H120").Value = '(L.) Soó'
$ws.Range("Q120").Value = 509705
$ws.Range("R120").Value = 6718923
$ws.Range("AC120").Value = 'Måttlig förekomst . inventering åt vasa vind'

# Row 127
$ws.Range("A127").Value = 130964546
$ws.Range("B127").Value = 92503
$ws.Range("D127").Value = 'VU'
$ws.Range("E127").Value = 898
$ws.Range("F127").Value = 'Blackticka'
$ws.Range("G127").Value = 'Steccherinum collabens'
$ws.Range("H127").Value = '(Fr.) Vesterholt'
$ws.Range("Q127").Value = 509515
$ws.Range("R127").Value = 6718886
$ws.Range("AC127").Value = 'Betydande förekomst . inventering åt vasa vind'

# Row 128
$ws.Range("A128").Value = 130964538
$ws.Range("B128").Value = 79243
$ws.Range("D128").Value = 'NT'
$ws.Range("E128").Value = 6425
$ws.Range("F128").Value = 'Garnlav'
$ws.Range("G128").Value = 'Alectoria sarmentosa'
$ws.Range("H128").Value = '(Ach.) Ach.'
$ws.Range("Q128").Value = 509875
$ws.Range("R128").Value = 6719025
$ws.Range("AC128").Value = 'Enstaka . inventering åt vasa vind'

# AX column normalization
$ws.Range("AX96").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'
$ws.Range("AX97").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'
$ws.Range("AX110").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'
$ws.Range("AX116").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'
$ws.Range("AX124").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'
$ws.Range("AX126").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'
